$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Clean up of parameters: update the "p0" (column D) start values for a
# handful of rate constants ahead of the next optimisation run.
$ws.Range("D7").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("D31").Value = 1

# Leave the sheet scrolled near the top with the cursor parked on the
# cell that was just edited (matches the saved view state).
$ws.Range("D7").Select()
